# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (columns A-J)
#   *_new -> *_FV2410   (columns L-U)
# Then wrap the data range in an Excel Table (ListObject) with autofilter,
# and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) -------------------------------------------------
$fv2404Headers = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

# Columns A-J (1-10) -> _FV2404
for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2404Headers[$i]
}

# Column K (11) stays "diff"

# Columns L-U (12-21) -> _FV2410
for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2410Headers[$i]
}

# --- Freeze header row (row 1) -------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Convert the data range into a Table (ListObject) --------------------------
$range = $ws.Range("A1:U77")
$list = $ws.ListObjects.Add(1, $range, 0, 1)
$list.Name = "Table1"

# Restore the original selection/cursor position.
[void]$ws.Range("A1").Select()
